# Updated symbol list on Thu Feb  2 11:21:05 UTC 2023 with GitHub Actions
# Applies the latest coin-price / 1h-volume snapshot to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new (text) value. Values are written with an
# explicit text NumberFormat so Excel keeps them as strings (matching the
# original inlineStr cell contents) instead of auto-converting numeric-
# looking text ("330.45", "7.34%", ...) into real numbers.
$updates = @{
    "D2" = "330.45"
    "E2" = "7.34%"
    "D3" = "40.13"
    "E3" = "8.14%"
    "D4" = "5.374"
    "E4" = "5.07%"
    "D5" = "0.08106"
    "E5" = "3.54%"
    "E6" = "3.00%"
    "D7" = "8.650"
    "E7" = "4.81%"
    "D8" = "1.929"
    "E8" = "1.67%"
    "E9" = "-0.35%"
    "D10" = "0.9426"
    "E10" = "2.37%"
    "D11" = "0.1358"
    "E11" = "25.28%"
    "D12" = "0.1969"
    "E12" = "3.91%"
    "D13" = "0.09336"
    "E13" = "5.44%"
    "D14" = "0.03561"
    "E14" = "6.28%"
    "D15" = "0.09595"
    "E15" = "0.18%"
    "D16" = "0.001326"
    "E16" = "-3.74%"
    "D18" = "3.366"
    "E18" = "-1.25%"
    "D19" = "0.3520"
    "D20" = "7.247"
    "E20" = "15.22%"
    "D21" = "0.1331"
    "E21" = "2.02%"
    "D22" = "0.2560"
    "E22" = "4.39%"
    "D23" = "0.04435"
    "E23" = "2.09%"
    "E24" = "2.42%"
    "D25" = "0.004275"
    "E25" = "0.30%"
    "E26" = "-14.27%"
    "D27" = "0.0003990"
    "E27" = "-0.06%"
    "D39" = "0.02491"
    "E39" = "14.79%"
    "D40" = "0.05222"
    "E40" = "4.03%"
    "D41" = "0.007562"
    "E41" = "0.17%"
    "D42" = "0.1429"
    "E42" = "5.98%"
    "E43" = "4.95%"
    "D44" = "0.002170"
    "E44" = "6.31%"
    "D45" = "0.01087"
    "E45" = "24.45%"
    "D46" = "0.00006564"
    "E46" = "0.47%"
    "E47" = "-0.07%"
    "D48" = "0.002400"
    "E48" = "139.46%"
    "E49" = "1.48%"
    "D50" = "0.00002100"
    "E50" = "-0.07%"
    "D51" = "0.0002000"
    "E51" = "-0.07%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}

Write-Output ("Updated " + $updates.Count + " cells")
